# Fruta / hortaliza, semanal
#
# Two new daily price records are inserted at the top of the Kiwi price
# block (rows 217-295), pushing the existing 79 records down by two rows
# (to 219-297). The sheet's used range grows from A1:T295 to A1:T297.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new records by inserting two blank rows at the
# current row 217 - everything below (217-295) shifts down to 219-297.
$ws.Rows.Item(217).Insert()
$ws.Rows.Item(217).Insert()

# --- New row 217 ---------------------------------------------------------
$ws.Range("A217").Value = 7
$ws.Range("B217").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C217").Value = "Ñuble"
$ws.Range("D217").Value2 = 45093
$ws.Range("E217").Value = 16
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100101
$ws.Range("H217").Value = "Berries"
$ws.Range("I217").Value = 100101007
$ws.Range("J217").Value = "Kiwi"
$ws.Range("K217").Value = "Hayward"
$ws.Range("L217").Value = "Especial"
$ws.Range("M217").Value = 130
$ws.Range("N217").Value = 12000
$ws.Range("O217").Value = 13000
$ws.Range("P217").Value = 12615
$ws.Range("Q217").Value = "$/bandeja 18 kilos"
$ws.Range("R217").Value = "Región de O'Higgins"
$ws.Range("S217").Value = 701
$ws.Range("T217").Value = 18

# --- New row 218 ---------------------------------------------------------
$ws.Range("A218").Value = 7
$ws.Range("B218").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C218").Value = "Ñuble"
$ws.Range("D218").Value2 = 45093
$ws.Range("E218").Value = 16
$ws.Range("F218").Value = "Fruta"
$ws.Range("G218").Value = 100101
$ws.Range("H218").Value = "Berries"
$ws.Range("I218").Value = 100101007
$ws.Range("J218").Value = "Kiwi"
$ws.Range("K218").Value = "Hayward"
$ws.Range("L218").Value = "Primera"
$ws.Range("M218").Value = 80
$ws.Range("N218").Value = 10000
$ws.Range("O218").Value = 10000
$ws.Range("P218").Value = 10000
$ws.Range("Q218").Value = "$/bandeja 18 kilos"
$ws.Range("R218").Value = "Región de O'Higgins"
$ws.Range("S218").Value = 556
$ws.Range("T218").Value = 18
